$wb = $excel.ActiveWorkbook

# --- Overview sheet: status for the 63769a2d... row (row 3) changes for both
# zh-cn (col B) and de-de (col C) from "Ready for handoff" to
# "Handback transform failed".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

# --- zh-cn sheet: Status column (C) for row 3 gets the same new status, and
# a new Error Detail (column L) value is added for that row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("L3").Value = "Handback file name: ogogteya.a4n is different with handoff file name: 63769a2d-717e-43e8-8343-4bbc1b4df3a9.f479f652c1648b498cd41734be3ea65c37dafe90.zh-cn."

# --- de-de sheet: same pair of edits as zh-cn, but with the de-de filename.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("L3").Value = "Handback file name: ogogteya.a4n is different with handoff file name: 63769a2d-717e-43e8-8343-4bbc1b4df3a9.f479f652c1648b498cd41734be3ea65c37dafe90.de-de."
